# Update retention metrics for ADD (2023 cohort, period_index 2, row 19)
# and for 2025 cohort (row 22): num_customers and cohort_size increase by 1,
# and retention_rate is recalculated as num_customers / cohort_size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: num_customers 59 -> 60; retention_rate recalculated (cohort_size stays 124)
$ws.Range("C19").Value = 60
$ws.Range("E19").Value = 60 / 124

# Row 22: num_customers 54 -> 55; cohort_size 54 -> 55; retention_rate stays 1
$ws.Range("C22").Value = 55
$ws.Range("D22").Value = 55
$ws.Range("E22").Value = 55 / 55
